$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove all existing rows/content so the sheet is rebuilt cleanly
$ws.Rows("1:26").Delete()

$ws.Range('B1').Value = 'Ementa atual:'
$ws.Range('C1').Value = 'Ementa modificada (dados modificados em vermelho):'

$ws.Range('B2').Value = 'LOM3108'
$ws.Range('C2').Value = 'LOM3108'

$ws.Range('A3').Value = 'Nome:'
$ws.Range('B3').Value = ' Projeto Integrado em Engenharia de Materiais II'
$ws.Range('C3').Value = ' Projeto Integrado em Engenharia de Materiais II'

$ws.Range('A4').Value = 'Name:'
$ws.Range('B4').Value = 'Materials Engineering Project II'
$ws.Range('C4').Value = 'Materials Engineering Project II'

$ws.Range('A5').Value = 'Créditos-aula:'
$ws.Range('B5').Value = '2'
$ws.Range('C5').Value = '2'

$ws.Range('A6').Value = 'Créditos-trabalho'
$ws.Range('B6').Value = '1'
$ws.Range('C6').Value = '1'

$ws.Range('A7').Value = 'Carga horária:'
$ws.Range('B7').Value = '60 h'
$ws.Range('C7').Value = '60 h'

$ws.Range('A8').Value = 'Ativação:'
$ws.Range('B8').Value = '01/01/2022'
$ws.Range('C8').Value = '01/01/2022'

$ws.Range('A9').Value = 'Semestre ideal:'
$ws.Range('B9').Value = 'EM-4'
$ws.Range('C9').Value = 'EM-4'

$ws.Range('A10').Value = 'Objetivos:'
$ws.Range('B10').Value = '7459752 - Maria Ismenia Sodero Toledo Faria'
$ws.Range('C10').Value = '7459752 - Maria Ismenia Sodero Toledo Faria'
$ws.Rows.Item(10).RowHeight = 60

$ws.Range('A11').Value = 'Objectives:'
$ws.Rows.Item(11).RowHeight = 60

$ws.Range('A12').Value = 'Docentes responsáveis:'

$ws.Range('A13').Value = 'Programa resumido:'
$ws.Range('B13').Value = '7459752 - Maria Ismenia Sodero Toledo Faria'
$ws.Range('C13').Value = '7459752 - Maria Ismenia Sodero Toledo Faria'
$ws.Rows.Item(13).RowHeight = 60

$ws.Range('A14').Value = 'Short syllabus:'
$ws.Rows.Item(14).RowHeight = 60

$ws.Range('A15').Value = 'Programa:'
$ws.Range('B15').Value = '2166002 - Sandra Giacomin Schneider'
$ws.Range('C15').Value = '2166002 - Sandra Giacomin Schneider'
$ws.Rows.Item(15).RowHeight = 120

$ws.Range('A16').Value = 'Syllabus:'
$ws.Rows.Item(16).RowHeight = 120

$ws.Range('A17').Value = 'Avaliação:'

$ws.Range('A18').Value = 'Método:'
$ws.Range('B18').Value = '1922320 - Sebastiao Ribeiro'
$ws.Range('C18').Value = '1922320 - Sebastiao Ribeiro'
$ws.Rows.Item(18).RowHeight = 60

$ws.Range('A19').Value = 'Critério:'
$ws.Range('B19').Value = 'O método utilizado tem por fundamento a aprendizagem baseada em projetos que visa desenvolver as competências técnicas relativas ao tema do projeto, bem como competências transversais, tais como: aprender a aprender, trabalho em equipe, relacionamento interpessoal, capacidade de comunicação oral e verbal e aspectos de liderança, dentre outros.Os alunos serão divididos em grupos que desenvolverão um projeto durante o semestre relacionado a um tema de Engenharia de Materiais, similar ao que eles irão encontrar na vida real no efetivo exercício de sua profissão.Cada grupo deverá buscar o conhecimento prático necessário para ser aplicado no desenvolvimento do projeto.As aulas ocorrerão por meio de uma reunião da equipe de trabalho para tratar do projeto; palestras e dinâmicas relativas ao tema do projeto, conduzidas por professores ou profissionais de empresas.'
$ws.Range('C19').Value = 'O método utilizado tem por fundamento a aprendizagem baseada em projetos que visa desenvolver as competências técnicas relativas ao tema do projeto, bem como competências transversais, tais como: aprender a aprender, trabalho em equipe, relacionamento interpessoal, capacidade de comunicação oral e verbal e aspectos de liderança, dentre outros.Os alunos serão divididos em grupos que desenvolverão um projeto durante o semestre relacionado a um tema de Engenharia de Materiais, similar ao que eles irão encontrar na vida real no efetivo exercício de sua profissão.Cada grupo deverá buscar o conhecimento prático necessário para ser aplicado no desenvolvimento do projeto.As aulas ocorrerão por meio de uma reunião da equipe de trabalho para tratar do projeto; palestras e dinâmicas relativas ao tema do projeto, conduzidas por professores ou profissionais de empresas.'
$ws.Rows.Item(19).RowHeight = 60

$ws.Range('A20').Value = 'Norma de recuperação:'
$ws.Range('B20').Value = 'A nota será individual e será a média ponderada de entregas do projeto, tais como: projeto preliminar, projeto final, envolvimento do aluno com o projeto, avaliação dos pares, autoavaliação e apresentação de trabalhos, dentre outros.O detalhamento dos pesos para ponderação da média da disciplina será definido por uma equipe de professores que atuarão na avaliação da disciplina.'
$ws.Range('C20').Value = 'A nota será individual e será a média ponderada de entregas do projeto, tais como: projeto preliminar, projeto final, envolvimento do aluno com o projeto, avaliação dos pares, autoavaliação e apresentação de trabalhos, dentre outros.O detalhamento dos pesos para ponderação da média da disciplina será definido por uma equipe de professores que atuarão na avaliação da disciplina.'
$ws.Rows.Item(20).RowHeight = 60

$ws.Range('A21').Value = 'Bibliografia:'
$ws.Range('B21').Value = 'não há'
$ws.Range('C21').Value = 'não há'
$ws.Rows.Item(21).RowHeight = 120

$ws.Range('A22').Value = 'Requisitos:'

$ws.Range('B23').Value = "LOM3104 -  Projeto Integrado em Engenharia de Materiais I  (Requisito fraco)`n"
$ws.Range('C23').Value = "LOM3104 -  Projeto Integrado em Engenharia de Materiais I  (Requisito fraco)`n"
$ws.Rows.Item(23).RowHeight = 30
